$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 'Multi-Utilities(18)'
$ws.Range("B2").Value = 0.4701696426296317
$ws.Range("A3").Value = 'Energy Equipment & Services(32)'
$ws.Range("B3").Value = 0.406557305300907
$ws.Range("A4").Value = 'Road & Rail(22)'
$ws.Range("B4").Value = 0.3673881964468256
$ws.Range("A5").Value = 'Banks(246)'
$ws.Range("B5").Value = 0.3368748846693499
$ws.Range("A6").Value = 'Electric Utilities(28)'
$ws.Range("B6").Value = 0.3205586396304085
$ws.Range("A7").Value = 'Marine(15)'
$ws.Range("B7").Value = 0.3070128247586821
$ws.Range("A8").Value = 'Building Products(23)'
$ws.Range("B8").Value = 0.302143338175924
$ws.Range("A9").Value = 'Auto Components(21)'
$ws.Range("B9").Value = 0.2938848054147093
$ws.Range("A10").Value = 'Machinery(85)'
$ws.Range("B10").Value = 0.2664101813557964
$ws.Range("A11").Value = 'Construction & Engineering(20)'
$ws.Range("B11").Value = 0.2616263087945421
$ws.Range("A12").Value = 'Trading Companies & Distributors(25)'
$ws.Range("B12").Value = 0.2552831537526278
$ws.Range("A13").Value = 'Specialty Retail(58)'
$ws.Range("B13").Value = 0.2498051986964366
$ws.Range("A14").Value = 'Thrifts & Mortgage Finance(47)'
$ws.Range("B14").Value = 0.2357986428000619
$ws.Range("A15").Value = 'Textiles, Apparel & Luxury Goods(29)'
$ws.Range("B15").Value = 0.2332191071178296
$ws.Range("A16").Value = 'Capital Markets(75)'
$ws.Range("B16").Value = 0.1909815052457282
$ws.Range("A17").Value = 'Hotels, Restaurants & Leisure(50)'
$ws.Range("B17").Value = 0.1870392878746558
$ws.Range("A18").Value = 'Insurance(75)'
$ws.Range("B18").Value = 0.1840892825553289
$ws.Range("A19").Value = 'Oil, Gas & Consumable Fuels(122)'
$ws.Range("B19").Value = 0.1837975971606767
$ws.Range("A20").Value = 'Semiconductors & Semiconductor Equipment(68)'
$ws.Range("B20").Value = 0.1803422401853819
$ws.Range("A21").Value = 'Professional Services(35)'
$ws.Range("B21").Value = 0.1672019193938927
$ws.Range("A22").Value = 'IT Services(52)'
$ws.Range("B22").Value = 0.1491692899065586
$ws.Range("A23").Value = 'Chemicals(51)'
$ws.Range("B23").Value = 0.1458712031784943
$ws.Range("A24").Value = 'Metals & Mining(89)'
$ws.Range("B24").Value = 0.1440840363727988
$ws.Range("A25").Value = 'Household Durables(39)'
$ws.Range("B25").Value = 0.1340439084648633
$ws.Range("A26").Value = 'Aerospace & Defense(37)'
$ws.Range("B26").Value = 0.1328854393022718
$ws.Range("A27").Value = 'Health Care Equipment & Supplies(83)'
$ws.Range("B27").Value = 0.1205744696534682
$ws.Range("A28").Value = 'Health Care Providers & Services(46)'
$ws.Range("B28").Value = 0.1164692854126994
$ws.Range("A29").Value = 'Commercial Services & Supplies(52)'
$ws.Range("B29").Value = 0.1076821198739539
$ws.Range("A30").Value = 'Software(66)'
$ws.Range("B30").Value = 0.08182842459818392
$ws.Range("A31").Value = 'Biotechnology(126)'
$ws.Range("B31").Value = 0.07576123345670216
